$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1039.2142
$ws.Range("I2").Value = 1099.9166
$ws.Range("J2").Value = 675
$ws.Range("K2").Value = 1099.9166
$ws.Range("L2").Value = 675
$ws.Range("M2").Value = -986.9166
$ws.Range("N2").Value = -901

$ws.Range("I5").Value = 203.72728
$ws.Range("J5").Value = 344.25
$ws.Range("K5").Value = 203.72728
$ws.Range("L5").Value = 344.25
$ws.Range("M5").Value = -88.72728000000001
$ws.Range("N5").Value = -574.25

$ws.Range("H58").Value = 7473.778
$ws.Range("I58").Value = 2038.3334
$ws.Range("K58").Value = 6115.0002
$ws.Range("M58").Value = -5965.0002

$ws.Range("H98").Value = 11303.733
$ws.Range("I98").Value = 13404.917
$ws.Range("J98").Value = 2899
$ws.Range("K98").Value = 13404.917
$ws.Range("L98").Value = 2899
$ws.Range("M98").Value = -11906.917
$ws.Range("N98").Value = -5895

$ws.Range("H100").Value = 2822.0908
$ws.Range("I100").Value = 1778.8572
$ws.Range("K100").Value = 1778.8572
$ws.Range("M100").Value = -1237.8572

$ws.Range("H106").Value = 3300
$ws.Range("I106").Value = 3214.2856
$ws.Range("K106").Value = 3214.2856
$ws.Range("M106").Value = -2583.2856

$ws.Range("H122").Value = 11303.733
$ws.Range("I122").Value = 13404.917
$ws.Range("J122").Value = 2899
$ws.Range("K122").Value = 40214.751
$ws.Range("L122").Value = 8697
$ws.Range("M122").Value = -37764.751
$ws.Range("N122").Value = -13597

$ws.Range("H137").Value = 4171.1284
$ws.Range("I137").Value = 2395.7188
$ws.Range("J137").Value = 12287.286
$ws.Range("K137").Value = 7187.1564
$ws.Range("L137").Value = 36861.858
$ws.Range("M137").Value = -4637.1564
$ws.Range("N137").Value = -41961.858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 39314.77
$ws.Range("I32").Value = 39129.773
$ws.Range("K32").Value = 39129.773
$ws.Range("M32").Value = -38842.773

$ws.Range("H61").Value = 16673772
$ws.Range("I61").Value = 27781822
$ws.Range("J61").Value = 11695.375
$ws.Range("K61").Value = 27781822
$ws.Range("L61").Value = 11695.375
$ws.Range("M61").Value = -27781610
$ws.Range("N61").Value = -12119.375

$ws.Range("H88").Value = 913.75
$ws.Range("J88").Value = 899
$ws.Range("L88").Value = 899
$ws.Range("N88").Value = -1711

$ws.Range("H91").Value = 913.75
$ws.Range("J91").Value = 899
$ws.Range("L91").Value = 899
$ws.Range("N91").Value = -3707

$ws.Range("H110").Value = 16669449
$ws.Range("I110").Value = 31251716
$ws.Range("K110").Value = 31251716
$ws.Range("M110").Value = -31249671

$ws.Range("H132").Value = 7150168
$ws.Range("I132").Value = 15390165
$ws.Range("J132").Value = 8837.532999999999
$ws.Range("K132").Value = 46170495
$ws.Range("L132").Value = 26512.599
$ws.Range("M132").Value = -46167965
$ws.Range("N132").Value = -31572.599

$ws.Range("H136").Value = 16673772
$ws.Range("I136").Value = 27781822
$ws.Range("J136").Value = 11695.375
$ws.Range("K136").Value = 83345466
$ws.Range("L136").Value = 35086.125
$ws.Range("M136").Value = -83342916
$ws.Range("N136").Value = -40186.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 85021.71000000001
$ws.Range("I86").Value = 1633.2307
$ws.Range("J86").Value = 183571.73
$ws.Range("K86").Value = 1633.2307
$ws.Range("L86").Value = 183571.73
$ws.Range("M86").Value = -510.2307000000001
$ws.Range("N86").Value = -185817.73

$ws.Range("H89").Value = 85021.71000000001
$ws.Range("I89").Value = 1633.2307
$ws.Range("J89").Value = 183571.73
$ws.Range("K89").Value = 8166.1535
$ws.Range("L89").Value = 917858.65
$ws.Range("M89").Value = -2550.1535
$ws.Range("N89").Value = -929090.65

$ws.Range("H94").Value = 1994.8636
$ws.Range("I94").Value = 2309.4
$ws.Range("J94").Value = 1732.75
$ws.Range("K94").Value = 2309.4
$ws.Range("L94").Value = 1732.75
$ws.Range("M94").Value = -1858.4
$ws.Range("N94").Value = -2634.75

$ws.Range("H105").Value = 52633316
$ws.Range("I105").Value = 62501376
$ws.Range("K105").Value = 62501376
$ws.Range("M105").Value = -62499629

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 497.1
$ws.Range("J22").Value = 596.6667
$ws.Range("L22").Value = 596.6667
$ws.Range("N22").Value = -1296.6667

$ws.Range("H31").Value = 5350.577
$ws.Range("I31").Value = 3005
$ws.Range("J31").Value = 6816.5625
$ws.Range("K31").Value = 3005
$ws.Range("L31").Value = 6816.5625
$ws.Range("M31").Value = -2710
$ws.Range("N31").Value = -7406.5625

$ws.Range("H34").Value = 5350.577
$ws.Range("I34").Value = 3005
$ws.Range("J34").Value = 6816.5625
$ws.Range("K34").Value = 3005
$ws.Range("L34").Value = 6816.5625
$ws.Range("M34").Value = -2803
$ws.Range("N34").Value = -7220.5625

$ws.Range("H94").Value = 2945.389
$ws.Range("I94").Value = 5959.143
$ws.Range("J94").Value = 1027.5454
$ws.Range("K94").Value = 5959.143
$ws.Range("L94").Value = 1027.5454
$ws.Range("M94").Value = -5508.143
$ws.Range("N94").Value = -1929.5454

$ws.Range("H99").Value = 6166.3335
$ws.Range("J99").Value = 6750
$ws.Range("L99").Value = 6750
$ws.Range("N99").Value = -9746

$ws.Range("H122").Value = 2215.4736
$ws.Range("I122").Value = 684.2308
$ws.Range("K122").Value = 2052.6924
$ws.Range("M122").Value = 397.3076000000001

$ws.Range("H126").Value = 6166.3335
$ws.Range("J126").Value = 6750
$ws.Range("L126").Value = 20250
$ws.Range("N126").Value = -25190

$ws.Range("H132").Value = 18934.611
$ws.Range("I132").Value = 1919.3704
$ws.Range("J132").Value = 69980.336
$ws.Range("K132").Value = 5758.1112
$ws.Range("L132").Value = 209941.008
$ws.Range("M132").Value = -3228.1112
$ws.Range("N132").Value = -215001.008

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 1100
$ws.Range("I8").Value = 1100
$ws.Range("K8").Value = 3300
$ws.Range("M8").Value = -3161

$ws.Range("H33").Value = 232.13333
$ws.Range("I33").Value = 135.6
$ws.Range("K33").Value = 813.5999999999999
$ws.Range("M33").Value = -530.5999999999999

$ws.Range("H37").Value = 98265.37
$ws.Range("J37").Value = 98265.37
$ws.Range("L37").Value = 294796.11
$ws.Range("N37").Value = -295020.11

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1138.8
$ws.Range("I97").Value = 1098.5
$ws.Range("J97").Value = 1300
$ws.Range("K97").Value = 1098.5
$ws.Range("L97").Value = 1300
$ws.Range("M97").Value = -602.5
$ws.Range("N97").Value = -2292

$ws.Range("H102").Value = 2735.75
$ws.Range("I102").Value = 2698.1428
$ws.Range("K102").Value = 2698.1428
$ws.Range("M102").Value = -1076.1428

$ws.Range("H113").Value = 1576
$ws.Range("I113").Value = 1495
$ws.Range("K113").Value = 1495
$ws.Range("M113").Value = 675

$ws.Range("H123").Value = 63164.332
$ws.Range("J123").Value = 63164.332
$ws.Range("L123").Value = 63164.332
$ws.Range("N123").Value = -68064.33199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 252251.25
$ws.Range("I7").Value = 335666.66
$ws.Range("J7").Value = 2005
$ws.Range("K7").Value = 335666.66
$ws.Range("L7").Value = 2005
$ws.Range("M7").Value = -335554.66
$ws.Range("N7").Value = -2229

$ws.Range("H61").Value = 3012.077
$ws.Range("I61").Value = 2215.8
$ws.Range("K61").Value = 2215.8
$ws.Range("M61").Value = -2013.8

$ws.Range("H113").Value = 3012.077
$ws.Range("I113").Value = 2215.8
$ws.Range("K113").Value = 2215.8
$ws.Range("M113").Value = -45.80000000000018

$ws.Range("H122").Value = 3297.5715
$ws.Range("I122").Value = 3297.5715
$ws.Range("K122").Value = 9892.7145
$ws.Range("M122").Value = -7442.7145

$ws.Range("H126").Value = 252251.25
$ws.Range("I126").Value = 335666.66
$ws.Range("J126").Value = 2005
$ws.Range("K126").Value = 1006999.98
$ws.Range("L126").Value = 6015
$ws.Range("M126").Value = -1004529.98
$ws.Range("N126").Value = -10955

$ws.Range("H132").Value = 4067.0444
$ws.Range("I132").Value = 2720.963
$ws.Range("J132").Value = 6086.1665
$ws.Range("K132").Value = 8162.889000000001
$ws.Range("L132").Value = 18258.4995
$ws.Range("M132").Value = -5632.889000000001
$ws.Range("N132").Value = -23318.4995

$ws.Range("H136").Value = 1313096.5
$ws.Range("I136").Value = 1555792.4
$ws.Range("K136").Value = 4667377.199999999
$ws.Range("M136").Value = -4664827.199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5058.905
$ws.Range("I81").Value = 4907.722
$ws.Range("J81").Value = 5966
$ws.Range("K81").Value = 9815.444
$ws.Range("L81").Value = 11932
$ws.Range("M81").Value = -8754.444
$ws.Range("N81").Value = -14054

$ws.Range("H84").Value = 5058.905
$ws.Range("I84").Value = 4907.722
$ws.Range("J84").Value = 5966
$ws.Range("K84").Value = 49077.22
$ws.Range("L84").Value = 59660
$ws.Range("M84").Value = -43773.22
$ws.Range("N84").Value = -70268

$ws.Range("H94").Value = 45948
$ws.Range("J94").Value = 45948
$ws.Range("L94").Value = 45948
$ws.Range("N94").Value = -47750

$ws.Range("H122").Value = 2242.5813
$ws.Range("I122").Value = 2189.0322
$ws.Range("K122").Value = 6567.096600000001
$ws.Range("M122").Value = -4117.096600000001

$ws.Range("H132").Value = 4878.2163
$ws.Range("I132").Value = 3708.0688
$ws.Range("K132").Value = 11124.2064
$ws.Range("M132").Value = -8594.206399999999

$ws.Range("H136").Value = 2314.85
$ws.Range("I136").Value = 1243.8334
$ws.Range("J136").Value = 5527.9
$ws.Range("K136").Value = 3731.5002
$ws.Range("L136").Value = 16583.7
$ws.Range("M136").Value = -1181.5002
$ws.Range("N136").Value = -21683.7
